{"js": "// The document contains a single 5-column table. Every 4th row (0, 4, 8,\n// 12, 16) holds five \"two-digit \u00f7 one-digit\" division answers; the rows in\n// between are blank spacer rows. The edit replaces the 25 answer strings,\n// in document order, with newly generated ones (same \"a\u00f7b=c, d\" shape, new\n// numbers). Mapping is strictly positional, since some old answer strings\n// repeat (e.g. \"68\u00f74=17, 0\" appears twice but is replaced differently each\n// time), so we address cells by (row, column) rather than by text search.\n\nconst newAnswers = [\n  [\"51\u00f72=25, 1\", \"93\u00f73=31, 0\", \"35\u00f74=8, 3\", \"89\u00f79=9, 8\", \"77\u00f72=38, 1\"],\n  [\"90\u00f72=45, 0\", \"58\u00f72=29, 0\", \"70\u00f74=17, 2\", \"25\u00f75=5, 0\", \"91\u00f72=45, 1\"],\n  [\"43\u00f78=5, 3\", \"23\u00f72=11, 1\", \"69\u00f79=7, 6\", \"88\u00f79=9, 7\", \"84\u00f78=10, 4\"],\n  [\"14\u00f75=2, 4\", \"28\u00f75=5, 3\", \"14\u00f76=2, 2\", \"56\u00f76=9, 2\", \"52\u00f77=7, 3\"],\n  [\"10\u00f79=1, 1\", \"80\u00f73=26, 2\", \"27\u00f73=9, 0\", \"22\u00f74=5, 2\", \"35\u00f73=11, 2\"],\n];\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\n// Answer rows are every 4th row starting at 0 (row, +3 blank rows, row, ...).\nfor (let block = 0; block < newAnswers.length; block++) {\n  const rowIndex = block * 4;\n  const rowValues = newAnswers[block];\n  for (let col = 0; col < rowValues.length; col++) {\n    const cell = table.getCell(rowIndex, col);\n    cell.value = rowValues[col];\n  }\n}\n\nawait context.sync();\n", "ps1": "# The document contains a single 5-column table. Every 4th row (1, 5, 9, 13,\n# 17 in Word's 1-based row numbering) holds five \"two-digit / one-digit\"\n# division answers; the rows in between are blank spacer rows. The edit\n# replaces the 25 answer strings, in document order, with newly generated\n# ones (same \"a/b=c, d\" shape, new numbers). Mapping is strictly positional,\n# since some old answer strings repeat (e.g. \"68/4=17, 0\" appears twice but\n# is replaced differently each time), so we address cells by (row, column)\n# rather than by text search/replace.\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n$newAnswers = @(\n    @(\"51\u00f72=25, 1\", \"93\u00f73=31, 0\", \"35\u00f74=8, 3\", \"89\u00f79=9, 8\", \"77\u00f72=38, 1\"),\n    @(\"90\u00f72=45, 0\", \"58\u00f72=29, 0\", \"70\u00f74=17, 2\", \"25\u00f75=5, 0\", \"91\u00f72=45, 1\"),\n    @(\"43\u00f78=5, 3\", \"23\u00f72=11, 1\", \"69\u00f79=7, 6\", \"88\u00f79=9, 7\", \"84\u00f78=10, 4\"),\n    @(\"14\u00f75=2, 4\", \"28\u00f75=5, 3\", \"14\u00f76=2, 2\", \"56\u00f76=9, 2\", \"52\u00f77=7, 3\"),\n    @(\"10\u00f79=1, 1\", \"80\u00f73=26, 2\", \"27\u00f73=9, 0\", \"22\u00f74=5, 2\", \"35\u00f73=11, 2\")\n)\n\n# Answer rows are every 4th row starting at 1 (Word rows are 1-based).\nfor ($block = 0; $block -lt $newAnswers.Length; $block++) {\n    $rowIndex = $block * 4 + 1\n    $rowValues = $newAnswers[$block]\n    for ($col = 0; $col -lt $rowValues.Length; $col++) {\n        $t.Cell($rowIndex, $col + 1).Range.Text = $rowValues[$col]\n    }\n}\n"}
